$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()), (''selector'', ''passthrough''),
                (''model'',
                 AdaBoostClassifier(estimator=DecisionTreeClassifier(class_weight=''balanced'',
                                                                     criterion=''entropy'',
                                                                     max_depth=1,
                                                                     max_features=''log2'',
                                                                     min_samples_leaf=3,
                                                                     min_samples_split=4,
                                                                     random_state=42),
                                    n_estimators=5, random_state=42))])'
$ws.Range("B2").Value = 0.6476190476190476
$ws.Range("C2").Value = '{''scaler'': RobustScaler(), ''model__n_estimators'': 5, ''model__estimator__min_samples_split'': 4, ''model__estimator__min_samples_leaf'': 3, ''model__estimator__max_features'': ''log2'', ''model__estimator__max_depth'': 1, ''model__estimator__criterion'': ''entropy'', ''model__estimator__class_weight'': ''balanced''}'
$ws.Range("D2").Value = 0.4285714285714285
$ws.Range("E2").Value = '[1 0 0 1 0 0 1 1 0 1 0 0]'
$ws.Range("F2").Value = '[0 1 1 0 1 1 1 1 1 1 1 0]'
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0.9767619047619047
$ws.Range("I2").Value = 0.003727604351283464
$ws.Range("J2").Value = 0.5495238095238095
$ws.Range("K2").Value = 0.07485113344482791

$ws.Range("A3").Value = 'Pipeline(steps=[(''scaler'', StandardScaler()), (''selector'', ''passthrough''),
                (''model'',
                 AdaBoostClassifier(estimator=DecisionTreeClassifier(class_weight=''balanced'',
                                                                     max_depth=5,
                                                                     max_features=''log2'',
                                                                     min_samples_leaf=2,
                                                                     min_samples_split=4,
                                                                     random_state=42),
                                    n_estimators=10, random_state=42))])'
$ws.Range("B3").Value = 0.6190476190476191
$ws.Range("C3").Value = '{''scaler'': StandardScaler(), ''model__n_estimators'': 10, ''model__estimator__min_samples_split'': 4, ''model__estimator__min_samples_leaf'': 2, ''model__estimator__max_features'': ''log2'', ''model__estimator__max_depth'': 5, ''model__estimator__criterion'': ''gini'', ''model__estimator__class_weight'': ''balanced''}'
$ws.Range("D3").Value = 0.4285714285714285
$ws.Range("E3").Value = '[1 0 1 0 0 0 0 1 1 0 1 1]'
$ws.Range("F3").Value = '[0 1 0 1 1 0 1 1 1 1 0 1]'
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 0.9733809523809523
$ws.Range("I3").Value = 0.004983754460820096
$ws.Range("J3").Value = 0.5329523809523808
$ws.Range("K3").Value = 0.0745426910615209

$ws.Range("A4").Value = 'Pipeline(steps=[(''scaler'', RobustScaler()), (''selector'', ''passthrough''),
                (''model'',
                 AdaBoostClassifier(estimator=DecisionTreeClassifier(max_depth=4,
                                                                     max_features=''log2'',
                                                                     min_samples_leaf=5,
                                                                     min_samples_split=6,
                                                                     random_state=42),
                                    n_estimators=5, random_state=42))])'
$ws.Range("B4").Value = 0.6
$ws.Range("C4").Value = '{''scaler'': RobustScaler(), ''model__n_estimators'': 5, ''model__estimator__min_samples_split'': 6, ''model__estimator__min_samples_leaf'': 5, ''model__estimator__max_features'': ''log2'', ''model__estimator__max_depth'': 4, ''model__estimator__criterion'': ''gini'', ''model__estimator__class_weight'': None}'
$ws.Range("D4").Value = 0.5333333333333333
$ws.Range("E4").Value = '[1 0 1 1 1 1 0 1 0 1 0 1]'
$ws.Range("F4").Value = '[0 0 0 1 1 1 1 0 1 1 1 0]'
$ws.Range("G4").Value = 42
$ws.Range("H4").Value = 0.983952380952381
$ws.Range("I4").Value = 0.003776593048049497
$ws.Range("J4").Value = 0.4895238095238096
$ws.Range("K4").Value = 0.08853720127267714
